$wb = $excel.ActiveWorkbook

# --- "summary" sheet: new 5th test-set row (row 8) ---
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("B8").Value = "5_wm"
$wsSummary.Range("C8").Value = "cityofnewyork"

# --- "baseline" sheet: matching new row (row 9) with measured figures ---
$wsBaseline = $wb.Worksheets.Item("baseline")
$wsBaseline.Range("B9").Value = "5_wm"
$wsBaseline.Range("C9").Value = "cityofnewyork"
$wsBaseline.Range("D9").Value = 1
$wsBaseline.Range("E9").Value = 1
$wsBaseline.Range("G9").Value = 0.15
$wsBaseline.Range("H9").Value = 0.15

# --- selections / active tab: user ends up looking at the "baseline" sheet ---
$wsSummary.Range("B22").Select()
$wsBaseline.Activate()
$wsBaseline.Range("H10").Select()
